$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (2006 / 宮城県保険環境センター entry) - rows 8 and 9 shift up
$ws.Rows.Item(7).Delete()
